$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Wnt6"
$ws.Range("C2").Value = "Fzd7"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.4655723333333334
$ws.Range("H2").Value = 1.396717
$ws.Range("I2").Value = 0.5605021543775982
$ws.Range("J2").Value = 0.5605021543775982
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.553279333333334
$ws.Range("N2").Value = 7.659838000000001
$ws.Range("O2").Value = 0.1645043904057808
$ws.Range("P2").Value = 0.1645043904057808
$ws.Range("Q2").Value = 1.188736216871778
$ws.Range("R2").Value = 10.698625951846
$ws.Range("S2").Value = 0.09220506522701363
$ws.Range("T2").Value = 0.09220506522701362

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Wnt6"
$ws.Range("C3").Value = "Fzd7"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.4655723333333334
$ws.Range("H3").Value = 1.396717
$ws.Range("I3").Value = 0.5605021543775982
$ws.Range("J3").Value = 0.5605021543775982
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 8.058662
$ws.Range("N3").Value = 24.175986
$ws.Range("O3").Value = 0.5192088709172035
$ws.Range("P3").Value = 0.5192088709172035
$ws.Range("Q3").Value = 3.751890070884667
$ws.Range("R3").Value = 33.76701063796201
$ws.Range("S3").Value = 0.2910176907210528
$ws.Range("T3").Value = 0.2910176907210528

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Wnt6"
$ws.Range("C4").Value = "Fzd7"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.4655723333333334
$ws.Range("H4").Value = 1.396717
$ws.Range("I4").Value = 0.5605021543775982
$ws.Range("J4").Value = 0.5605021543775982
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 4.909099333333334
$ws.Range("N4").Value = 14.727298
$ws.Range("O4").Value = 0.3162867386770157
$ws.Range("P4").Value = 0.3162867386770157
$ws.Range("Q4").Value = 2.285540831185112
$ws.Range("R4").Value = 20.56986748066601
$ws.Range("S4").Value = 0.1772793984295317
$ws.Range("T4").Value = 0.1772793984295317

# Row 5
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Wnt6"
$ws.Range("C5").Value = "Fzd7"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.365062
$ws.Range("H5").Value = 1.095186
$ws.Range("I5").Value = 0.4394978456224017
$ws.Range("J5").Value = 0.4394978456224018
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.553279333333334
$ws.Range("N5").Value = 7.659838000000001
$ws.Range("O5").Value = 0.1645043904057808
$ws.Range("P5").Value = 0.1645043904057808
$ws.Range("Q5").Value = 0.9321052599853334
$ws.Range("R5").Value = 8.388947339868
$ws.Range("S5").Value = 0.07229932517876715
$ws.Range("T5").Value = 0.07229932517876715

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Wnt6"
$ws.Range("C6").Value = "Fzd7"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.365062
$ws.Range("H6").Value = 1.095186
$ws.Range("I6").Value = 0.4394978456224017
$ws.Range("J6").Value = 0.4394978456224018
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 8.058662
$ws.Range("N6").Value = 24.175986
$ws.Range("O6").Value = 0.5192088709172035
$ws.Range("P6").Value = 0.5192088709172035
$ws.Range("Q6").Value = 2.941911267044
$ws.Range("R6").Value = 26.477201403396
$ws.Range("S6").Value = 0.2281911801961506
$ws.Range("T6").Value = 0.2281911801961506

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Wnt6"
$ws.Range("C7").Value = "Fzd7"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.365062
$ws.Range("H7").Value = 1.095186
$ws.Range("I7").Value = 0.4394978456224017
$ws.Range("J7").Value = 0.4394978456224018
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 4.909099333333334
$ws.Range("N7").Value = 14.727298
$ws.Range("O7").Value = 0.3162867386770157
$ws.Range("P7").Value = 0.3162867386770157
$ws.Range("Q7").Value = 1.792125620825334
$ws.Range("R7").Value = 16.129130587428
$ws.Range("S7").Value = 0.139007340247484
$ws.Range("T7").Value = 0.139007340247484
